# Update "想去人数" (want-to-go count) figures in column F
# for both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 5598
$ws1.Range("F8").Value = 6464
$ws1.Range("F11").Value = 1420
$ws1.Range("F12").Value = 48
$ws1.Range("F13").Value = 46

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 5598
$ws4.Range("F9").Value = 6464
$ws4.Range("F12").Value = 1420
$ws4.Range("F13").Value = 48
$ws4.Range("F14").Value = 46
